$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.621.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +7.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.815.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.25%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.19'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.57%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3529'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.239'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07752'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +12.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.646'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.219'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.812.38'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.48%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06698'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '86.96'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +7.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9986'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.85'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +10.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.544'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.21'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.530.63'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.463'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +15.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.504'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.19'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.018.16'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '137.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.406'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.081'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.98'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08842'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.723'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.668'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7114'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +16.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06567'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02429'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2272'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.020'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.59%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6631'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +13.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9987'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.984'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.67%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.42'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07371'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.93'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.04%  '
